# Auto-generated edit script: refresh crypto price/volume data
# and swap Hedera / WEMIXToken row order (rows 34-35).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.920.82'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '2.535.42'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.93'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.80'
$ws.Range('E6').Value = '  +2.87%  '
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.79'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('E13').Value = '  +1.73%  '
$ws.Range('D14').Value = '2.924.60'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.50'
$ws.Range('E15').Value = '  -2.69%  '
$ws.Range('D16').Value = '2.558.54'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').Value = '42.882.52'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.69'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.39'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.80'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '243.94'
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.04'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.62'
$ws.Range('E27').Value = '  -4.27%  '
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.23'
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.88'
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.92'
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.84'
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('E33').Value = '  +7.41%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.67'
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0792'
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.32'
$ws.Range('E36').Value = '  -1.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.14'
$ws.Range('E37').Value = '  -4.71%  '
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('E41').Value = '  +3.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.79'
$ws.Range('E42').Value = '  -3.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.33'
$ws.Range('E43').Value = '  +4.49%  '
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0299'
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('D46').Value = '2.004.94'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.24'
$ws.Range('E47').Value = '  +3.19%  '
$ws.Range('D48').Value = '2.778.25'
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.192'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '79.93'
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.44'
$ws.Range('E51').Value = '  -0.93%  '
